$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1.xml) ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Rows 10 and 11 both held a duplicate "Contact" / "No display for ContactDetail" pair.
# Delete one of them so the sheet collapses back down to 20 rows (A1:B20), then overwrite
# the remaining row with the new "Jurisdiction" / "United States of America" pair.
$ws1.Rows.Item(11).Delete()

$ws1.Range("B3").Value = "6.0.0"
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws1.Range("B9").Value = "Alvearie Team"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# --- Sheet "Elements" (sheet2.xml) ---
$ws2 = $wb.Worksheets.Item("Elements")

# Row 2 (the root "Extension" element) gets its Short/Definition replaced with the
# submission-type-specific text instead of the generic Extension placeholder text.
$ws2.Range("K2").Value = "Submission Type"
$ws2.Range("L2").Value = "Customer-specific code for the type of electronic submission"
